# Update the cryptos symbol list with the latest scraped prices/volumes.
# Numeric-looking strings must stay TEXT (as in the source sheet), so we
# briefly force a text number-format before assigning the value and then
# restore the default "Normal" style so no visible formatting changes.

function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
Set-TextCell $ws "D2" "293.09"
Set-TextCell $ws "E2" "-3.18%"

# Row 3 - OKB
Set-TextCell $ws "D3" "30.94"
Set-TextCell $ws "E3" "-3.70%"

# Row 4 - HuobiToken
Set-TextCell $ws "D4" "4.860"
Set-TextCell $ws "E4" "-2.27%"

# Row 5 - Cronos
Set-TextCell $ws "D5" "0.07267"
Set-TextCell $ws "E5" "-8.22%"

# Row 6 - was FTXToken, now KuCoinToken
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell $ws "D6" "7.676"
Set-TextCell $ws "E6" "-2.34%"

# Row 7 - was KuCoinToken, now FTXToken
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell $ws "D7" "1.756"
Set-TextCell $ws "E7" "-17.57%"

# Row 8 - GateToken
Set-TextCell $ws "D8" "3.765"
Set-TextCell $ws "E8" "-1.12%"

# Row 9 - MXToken
Set-TextCell $ws "D9" "0.9022"
Set-TextCell $ws "E9" "-2.67%"

# Row 10 - WazirX
Set-TextCell $ws "D10" "0.1653"
Set-TextCell $ws "E10" "-6.08%"

# Row 11 - LiechtensteinCryptoassetsExchange
Set-TextCell $ws "D11" "0.07587"
Set-TextCell $ws "E11" "-5.57%"

# Row 12 - MandalaExchangeToken
Set-TextCell $ws "D12" "0.08177"
Set-TextCell $ws "E12" "-6.77%"

# Row 13 - BitrueCoin
Set-TextCell $ws "D13" "0.03037"
Set-TextCell $ws "E13" "-3.80%"

# Row 14 - BitMartToken
Set-TextCell $ws "D14" "0.1000"
Set-TextCell $ws "E14" "-0.44%"

# Row 15
Set-TextCell $ws "D15" "0.001506"
Set-TextCell $ws "E15" "-0.65%"

# Row 16
Set-TextCell $ws "D16" "0.005654"
Set-TextCell $ws "E16" "-5.40%"

# Row 17
Set-TextCell $ws "D17" "3.458"
Set-TextCell $ws "E17" "-0.35%"

# Row 18
Set-TextCell $ws "D18" "2.106"
Set-TextCell $ws "E18" "-7.61%"

# Row 19
Set-TextCell $ws "D19" "0.3296"
Set-TextCell $ws "E19" "0.27%"

# Row 20
Set-TextCell $ws "D20" "0.1305"
Set-TextCell $ws "E20" "1.23%"

# Row 21
Set-TextCell $ws "D21" "4.363"
Set-TextCell $ws "E21" "3.93%"

# Row 22
Set-TextCell $ws "D22" "0.2004"
Set-TextCell $ws "E22" "11.97%"

# Row 23
Set-TextCell $ws "D23" "0.04493"
Set-TextCell $ws "E23" "-2.47%"

# Row 24 (only E changed)
Set-TextCell $ws "E24" "-1.67%"

# Row 25 (only D changed)
Set-TextCell $ws "D25" "0.004040"

# Row 26 (only E changed)
Set-TextCell $ws "E26" "0.19%"

# Row 39
Set-TextCell $ws "D39" "0.01650"
Set-TextCell $ws "E39" "-5.15%"

# Row 40
Set-TextCell $ws "D40" "0.04370"
Set-TextCell $ws "E40" "-9.08%"

# Row 41
Set-TextCell $ws "D41" "0.007409"
Set-TextCell $ws "E41" "0.86%"

# Row 42
Set-TextCell $ws "D42" "0.1320"
Set-TextCell $ws "E42" "-3.52%"

# Row 43 (only E changed)
Set-TextCell $ws "E43" "-12.75%"

# Row 44
Set-TextCell $ws "D44" "0.01026"
Set-TextCell $ws "E44" "-7.08%"

# Row 45
Set-TextCell $ws "D45" "0.00005669"
Set-TextCell $ws "E45" "-5.80%"

# Row 46
Set-TextCell $ws "D46" "0.00000000752"
Set-TextCell $ws "E46" "0.38%"

# Row 47
Set-TextCell $ws "D47" "2.174"
Set-TextCell $ws "E47" "164.96%"

# Row 48 (only E changed)
Set-TextCell $ws "E48" "-29.03%"

# Row 49
Set-TextCell $ws "D49" "0.00002105"
Set-TextCell $ws "E49" "0.38%"

# Row 50
Set-TextCell $ws "D50" "0.0002005"
Set-TextCell $ws "E50" "0.38%"
